$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: Approved/Rejected column now "Approved" (was "Rejected"),
#     and the ReasonToReject note that used to live on row 8 is gone.
$ws.Range("I8").Value = "Approved"
$ws.Range("J8").ClearContents()

# --- Row 10: Approved/Rejected column now "Rejected" (was "Approved"),
#     and the ReasonToReject note moved here.
$ws.Range("I10").Value = "Rejected"
$ws.Range("J10").Value = "Indetail explaination"

# --- View state: scroll the window so column H is left-most visible,
#     and leave the active selection on E8.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E8").Select()
